# Auto-generated edit script applying cryptos list price/volume updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.908.93"
$ws.Range("E2").Value = "  +6.89%  "
$ws.Range("D3").Value = "1.744.07"
$ws.Range("E3").Value = "  +5.48%  "
$ws.Range("D4").Value = "'1.004"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'228.17"
$ws.Range("E5").Value = "  +4.29%  "
$ws.Range("D6").Value = "'0.5462"
$ws.Range("E6").Value = "  +3.91%  "
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("D8").Value = "'0.2779"
$ws.Range("E8").Value = "  +4.16%  "
$ws.Range("D9").Value = "'0.06762"
$ws.Range("E9").Value = "  +6.20%  "
$ws.Range("D10").Value = "'21.89"
$ws.Range("E10").Value = "  +6.34%  "
$ws.Range("D11").Value = "'0.07797"
$ws.Range("E11").Value = "  +1.37%  "
$ws.Range("D12").Value = "'4.711"
$ws.Range("D13").Value = "1.764.00"
$ws.Range("E13").Value = "  +6.99%  "
$ws.Range("D14").Value = "1.985.45"
$ws.Range("E14").Value = "  +5.56%  "
$ws.Range("D15").Value = "'0.5992"
$ws.Range("E15").Value = "  +6.83%  "
$ws.Range("D16").Value = "0.0₅8417"
$ws.Range("E16").Value = "  +2.24%  "
$ws.Range("D17").Value = "'69.04"
$ws.Range("E17").Value = "  +5.53%  "
$ws.Range("D18").Value = "27.920.80"
$ws.Range("E18").Value = "  +6.92%  "
$ws.Range("D19").Value = "'224.73"
$ws.Range("E19").Value = "  +17.50%  "
$ws.Range("D20").Value = "'4.854"
$ws.Range("E20").Value = "  +3.20%  "
$ws.Range("D21").Value = "'1.002"
$ws.Range("E21").Value = "  -0.19%  "
$ws.Range("E22").Value = "  +5.43%  "
$ws.Range("D23").Value = "'6.254"
$ws.Range("E23").Value = "  +4.50%  "
$ws.Range("E24").Value = "  -0.16%  "
$ws.Range("D25").Value = "'146.32"
$ws.Range("E25").Value = "  +0.20%  "
$ws.Range("D26").Value = "'0.1254"
$ws.Range("E26").Value = "  +4.49%  "
$ws.Range("D27").Value = "'1.688"
$ws.Range("E27").Value = "  +12.70%  "
$ws.Range("D28").Value = "'7.480"
$ws.Range("E28").Value = "  +2.97%  "
$ws.Range("D29").Value = "'17.21"
$ws.Range("E29").Value = "  +7.87%  "
$ws.Range("D30").Value = "'0.05701"
$ws.Range("E30").Value = "  +0.77%  "
$ws.Range("D31").Value = "'1.318"
$ws.Range("E31").Value = "  +3.69%  "
$ws.Range("D32").Value = "'3.707"
$ws.Range("D33").Value = "'3.530"
$ws.Range("E33").Value = "  +4.16%  "
$ws.Range("D34").Value = "'1.693"
$ws.Range("E34").Value = "  +7.04%  "
$ws.Range("D35").Value = "'0.9782"
$ws.Range("E35").Value = "  +3.22%  "
$ws.Range("D36").Value = "'2.864"
$ws.Range("E36").Value = "  +2.30%  "
$ws.Range("D37").Value = "'2.448"
$ws.Range("E37").Value = "  +1.73%  "
$ws.Range("D38").Value = "'0.5989"
$ws.Range("E38").Value = "  +3.37%  "
$ws.Range("D39").Value = "'0.01670"
$ws.Range("D40").Value = "'5.966"
$ws.Range("E40").Value = "  -0.25%  "
$ws.Range("D41").Value = "'0.8518"
$ws.Range("E41").Value = "  +1.24%  "
$ws.Range("D42").Value = "1.049.24"
$ws.Range("E42").Value = "  +2.30%  "
$ws.Range("E43").Value = "  -0.05%  "
$ws.Range("D44").Value = "'102.28"
$ws.Range("D45").Value = "1.890.76"
$ws.Range("E45").Value = "  +5.51%  "
$ws.Range("D46").Value = "0.0₈117"
$ws.Range("E46").Value = "  +13.11%  "
$ws.Range("D47").Value = "'59.80"
$ws.Range("E47").Value = "  +2.11%  "
$ws.Range("D48").Value = "'8.330"
$ws.Range("E48").Value = "  +3.33%  "
$ws.Range("D49").Value = "'0.4438"
$ws.Range("E49").Value = "  +2.20%  "
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").Value = "'0.05336"
$ws.Range("E50").Value = "  -0.01%  "
$ws.Range("B51").Value = "Frax"
$ws.Range("C51").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D51").Value = "'1.003"
$ws.Range("E51").Value = "  -0.19%  "
